$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholder: cached datetimeFigureOut text bumped by one
#    day (2016/2/2 -> 2016/2/3) on the slide master and on every slide layout.
# ---------------------------------------------------------------------------
$targetOld = "2016/2/2"
$targetNew = "2016/2/3"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
            if ($isDatePlaceholder -and ($shp.TextFrame.TextRange.Text -eq $targetOld)) {
                $shp.TextFrame.TextRange.Text = $targetNew
            }
        }
    }
}

# Slide master
Update-DatePlaceholder($p.SlideMaster.Shapes)

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder($layout.Shapes)
}

# ---------------------------------------------------------------------------
# 2) Slide 2: fix "Ilinq" typo -> "ILinq" and split the following
#    "(Of Integer) " run into "(Of " / "Integer) " (as happens when the
#    "Integer" word is reselected and retyped).
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$grp = $slide2.Shapes.Item(1)
$codeShape = $grp.GroupItems.Item(10)
$codeTextRange = $codeShape.TextFrame.TextRange
$para = $codeTextRange.Paragraphs(3, 1)

# "Ilinq" -> "ILinq"
$ilinqRun = $para.Characters(15, 5)
$ilinqRun.Text = "ILinq"

# "(Of Integer) " run splits into "(Of " + "Integer) "
$integerPart = $para.Characters(24, 9)
$integerPart.Text = "Integer) "
